# "Fruta / hortaliza, semanal"
#
# A new weekly price-report record is inserted for the Espinaca
# (La Araucania / Vega Modelo de Temuco) data set. The new observation
# belongs between the existing row 54 (date 44172) and what used to be
# row 55 (date 44162), so it is inserted as the new row 55 - every
# following record (old rows 55-112) shifts down by one row to 56-113.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 55; Excel shifts rows 55-112 down to 56-113,
# carrying their values/styles with them (including the date style on D).
$ws.Rows.Item(55).Insert()

# Populate the newly inserted row 55 with the new weekly observation.
$ws.Range("A55").Value2 = 10
$ws.Range("B55").Value2 = "Vega Modelo de Temuco"
$ws.Range("C55").Value2 = "La Araucanía"
$ws.Range("D55").Value2 = 44546
$ws.Range("E55").Value2 = 9
$ws.Range("F55").Value2 = 100112012
$ws.Range("G55").Value2 = "Espinaca"
$ws.Range("H55").Value2 = "Sin especificar"
$ws.Range("I55").Value2 = "Primera"
$ws.Range("J55").Value2 = 115
$ws.Range("K55").Value2 = 8000
$ws.Range("L55").Value2 = 9000
$ws.Range("M55").Value2 = 8565
$ws.Range("N55").Value2 = "$/docena de atados"
$ws.Range("O55").Value2 = "Región de La Araucanía"
$ws.Range("P55").Value2 = 2855
$ws.Range("Q55").Value2 = 3
$ws.Range("R55").Value2 = "Hortaliza"
